$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy the existing header formatting
# (font/border/alignment) from G1 so the new header cell matches its
# siblings, then set the header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
